$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to the latest scraped values.
# The cells hold text (not numeric) data in the source workbook, so the
# number format is forced to Text before assignment and cleared again
# afterwards to avoid leaving a residual cell style behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '305.91'
Set-TextValue 'E2' '1.27%'
Set-TextValue 'D3' '36.28'
Set-TextValue 'E3' '-2.01%'
Set-TextValue 'D4' '5.054'
Set-TextValue 'E4' '0.94%'
Set-TextValue 'D5' '0.07922'
Set-TextValue 'E5' '3.11%'
Set-TextValue 'D6' '2.235'
Set-TextValue 'E6' '8.24%'
Set-TextValue 'D7' '8.001'
Set-TextValue 'E7' '0.45%'
Set-TextValue 'D8' '0.9272'
Set-TextValue 'E8' '1.15%'
Set-TextValue 'D9' '0.09833'
Set-TextValue 'E9' '2.99%'
Set-TextValue 'D10' '0.1872'
Set-TextValue 'E10' '1.14%'
Set-TextValue 'D11' '0.09009'
Set-TextValue 'E11' '4.84%'
Set-TextValue 'D12' '0.03736'
Set-TextValue 'E12' '3.53%'
Set-TextValue 'D13' '0.09916'
Set-TextValue 'E13' '-0.75%'
Set-TextValue 'D14' '0.001427'
Set-TextValue 'E14' '-3.48%'
Set-TextValue 'D15' '0.005634'
Set-TextValue 'E15' '-2.41%'
Set-TextValue 'D16' '3.466'
Set-TextValue 'E16' '-0.21%'
Set-TextValue 'D17' '4.152'
Set-TextValue 'E17' '2.80%'
Set-TextValue 'D18' '2.633'
Set-TextValue 'E18' '12.80%'
Set-TextValue 'D19' '0.3370'
Set-TextValue 'E19' '0.70%'
Set-TextValue 'D20' '0.1319'
Set-TextValue 'E20' '-1.16%'
Set-TextValue 'D21' '5.064'
Set-TextValue 'E21' '1.95%'
Set-TextValue 'D22' '0.2246'
Set-TextValue 'E22' '1.26%'
Set-TextValue 'D23' '0.04564'
Set-TextValue 'E23' '-1.00%'
Set-TextValue 'D24' '0.001237'
Set-TextValue 'E24' '-0.35%'
Set-TextValue 'D25' '0.004787'
Set-TextValue 'E25' '-5.87%'
Set-TextValue 'D26' '0.0001298'
Set-TextValue 'E26' '-8.03%'
Set-TextValue 'D39' '0.01918'
Set-TextValue 'E39' '10.06%'
Set-TextValue 'E40' '6.44%'
Set-TextValue 'D41' '0.007747'
Set-TextValue 'E41' '0.52%'
Set-TextValue 'D42' '0.1394'
Set-TextValue 'E42' '0.14%'
Set-TextValue 'D43' '0.007794'
Set-TextValue 'E43' '-2.48%'
Set-TextValue 'D44' '0.002145'
Set-TextValue 'E44' '-1.20%'
Set-TextValue 'E45' '15.55%'
Set-TextValue 'D46' '0.00006167'
Set-TextValue 'E46' '-2.09%'
Set-TextValue 'D47' '0.00000000749'
Set-TextValue 'E47' '-0.90%'
Set-TextValue 'D48' '51.76'
Set-TextValue 'E48' '51.19%'
Set-TextValue 'D49' '0.001798'
Set-TextValue 'E49' '-10.82%'
Set-TextValue 'D50' '0.00002098'
Set-TextValue 'E50' '-0.90%'
Set-TextValue 'D51' '0.0001998'
Set-TextValue 'E51' '-0.90%'
